$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the 2nd and 3rd data rows ("B" and "C" sub-rows) within every
# 4-row year block (A, B, C, D), for columns A through E.
$blockStarts = 2..68 | Where-Object { (($_ - 2) % 4) -eq 0 }

foreach ($start in $blockStarts) {
    $rowB = $start + 1
    $rowC = $start + 2

    for ($col = 1; $col -le 5; $col++) {
        $cellB = $ws.Cells.Item($rowB, $col)
        $cellC = $ws.Cells.Item($rowC, $col)

        $valB = $cellB.Value2
        $valC = $cellC.Value2

        # Cells that are already blank placeholders (empty inline strings,
        # e.g. column D for most rows) read back as an empty string;
        # writing that back would strip the cell's existing (empty-string)
        # representation, so only write back when there is an actual value.
        if (-not [string]::IsNullOrEmpty($valC)) { $cellB.Value2 = $valC }
        if (-not [string]::IsNullOrEmpty($valB)) { $cellC.Value2 = $valB }
    }
}

# Remove the now-unneeded "产销率" (F) and "销售量" (G) columns entirely.
$ws.Range("F1:G69").EntireColumn.Delete()
